# actualizacion feb 17 - 23:06
#
# Updates the "feb2025" sheet: fills in the february payments (column C =
# "pago1") with 65000 for every resident that paid, except the unoccupied
# house (Rafael, row 21) which paid a reduced 50000. Rows that already show
# a second payment ("pago2", column D) get that column updated too, mirroring
# the pattern already present on the other month sheets (e.g. ene2025).
# Also switches the active tab from "jun2025" back to "feb2025".

$wb = $excel.ActiveWorkbook

$wsFeb = $wb.Worksheets.Item("feb2025")

# pago1 (column C) - rows that got paid this period
$febPaidRows = 2,3,4,5,7,9,10,13,15,16,17,18,20,22,23,24
foreach ($r in $febPaidRows) {
    $wsFeb.Cells.Item($r, 3).Value = 65000
}

# row 21 = Rafael (casa desocupada) - reduced quota
$wsFeb.Cells.Item(21, 3).Value = 50000

# pago2 (column D) - rows where both installments are already settled
$wsFeb.Cells.Item(10, 4).Value = 65000
$wsFeb.Cells.Item(13, 4).Value = 65000
$wsFeb.Cells.Item(21, 4).Value = 50000
$wsFeb.Cells.Item(22, 4).Value = 65000

# Keep the shared "Rafael (casa desocupada)" label consistent across every
# month sheet that references it.
foreach ($name in @("feb2025","mar2025","abr2025","may2025","jun2025")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B21").Value = "Rafael (casa desocupada)"
}

# Move the selection on jun2025 (previously the active tab) off of F15...
$wsJun = $wb.Worksheets.Item("jun2025")
$wsJun.Range("E12").Select()

# ...and make feb2025 the active tab/selection instead.
$wsFeb.Activate()
$wsFeb.Range("D13").Select()
